$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths: A=36, B=5, C=5, D=5, E=13
$ws.Columns.Item(1).ColumnWidth = 35.083333333333336
$ws.Columns.Item(2).ColumnWidth = 4.083333333333333
$ws.Columns.Item(3).ColumnWidth = 4.083333333333333
$ws.Columns.Item(4).ColumnWidth = 4.083333333333333
$ws.Columns.Item(5).ColumnWidth = 12.083333333333332

# Replace the "-" placeholders in B2:E8 with numeric 0, formatted as "0.0"
$rng = $ws.Range("B2:E8")
$rng.Value = 0
$rng.NumberFormat = "0.0"
